# Applies the "Changes in handling of rows/columns" commit:
#   - Rename four worksheet tabs
#   - Move the active sheet/tab from "Contracts" to "APB Unit Cost History"
#     (renamed "Unit Cost History"), which also moves the saved selection
#     (tabSelected) and the workbook's activeTab index
#   - Update the saved cell selection on two sheets that scrolled/moved
#     their cursor as part of the same session

$wb = $excel.ActiveWorkbook

# --- Rename sheets -------------------------------------------------------
$wsFunding = $wb.Worksheets.Item("Cost and Funding")
$wsFunding.Name = "Funding Summary"

$wsFMS = $wb.Worksheets.Item("FMS")
$wsFMS.Name = "Foreign Military Sales"

$wsUnitCost = $wb.Worksheets.Item("Unit Cost Report")
$wsUnitCost.Name = "Unit Cost"

$wsUnitCostHistory = $wb.Worksheets.Item("APB Unit Cost History")
$wsUnitCostHistory.Name = "Unit Cost History"

# --- Update saved selections on sheets that were scrolled/navigated ------
$wsFunding.Activate() | Out-Null
$wsFunding.Range("A171").Select() | Out-Null

$wsUnitCost.Activate() | Out-Null
$wsUnitCost.Range("B30").Select() | Out-Null

# --- Finally land on "Unit Cost History" so it becomes the active tab ----
$wsUnitCostHistory.Activate() | Out-Null
